# improved version wit main loop and output fix
#
# Refreshes the two "Service Tag" summary blocks on the active sheet with a
# new scan's numbers:
#   - first block (rows 1-6):  service tag, IP, and the four pass/error counts
#   - second block (rows 8-13): IP + the four pass/error counts (tag unchanged)
#
# Several of the "numbers" are actually stored as plain text in the sheet
# (no number formatting is used anywhere), so when a replacement value looks
# numeric we enter it with a leading apostrophe (classic "force text" COM
# trick) and then snap the cell style back to Normal so we don't leave a
# lingering custom number format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    if ($newValue -match '^-?\d+(\.\d+)?$') {
        # numeric-looking -> force text entry so it isn't reinterpreted as a number
        $ws.Range($cellRef).Value = "'" + $newValue
        $ws.Range($cellRef).Style = "Normal"
    } else {
        $ws.Range($cellRef).Value = $newValue
    }
}

# row -> (ip/conf-passed/conf-error/hw-passed/hw-error or tag) updates, per block
$blocks = @(
    @{ Tag = "B1"; IP = "B2"; ConfPassed = "B3"; ConfError = "B4"; HwPassed = "B5"; HwError = "B6";
       TagValue = "C1VV2S2"; IPValue = "0.0.0.0"; ConfPassedValue = "1602"; ConfErrorValue = "0"; HwPassedValue = "70"; HwErrorValue = "1" },
    @{ Tag = $null;  IP = "B9"; ConfPassed = "B10"; ConfError = "B11"; HwPassed = "B12"; HwError = "B13";
       TagValue = $null; IPValue = "0.0.0.0"; ConfPassedValue = "1601"; ConfErrorValue = "1"; HwPassedValue = "71"; HwErrorValue = "0" }
)

foreach ($block in $blocks) {
    if ($block.Tag) {
        Set-TextValue $block.Tag $block.TagValue
    }
    Set-TextValue $block.IP $block.IPValue
    Set-TextValue $block.ConfPassed $block.ConfPassedValue
    Set-TextValue $block.ConfError $block.ConfErrorValue
    Set-TextValue $block.HwPassed $block.HwPassedValue
    Set-TextValue $block.HwError $block.HwErrorValue
}
